$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the header formatting (font + fill + etc.) from F1 into F2 so the
# new cell picks up the same font (white text) that the other header
# cells use, then give it its own distinct fill color.
$ws.Range("F1").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the new cell's text (goes into sharedStrings.xml automatically).
$ws.Range("F2").Value2 = "Need to Modularize error responses "

# Give F2 its own purple fill (RGB 0x70,0x30,0xA0 -> FF7030A0).
$ws.Range("F2").Interior.Color = 10498160

# Update the current selection on the sheet.
[void]$ws.Range("E9:E10").Select()
